$d = $word.ActiveDocument

# The document currently has:
#   ... "use " [bookmark _GoBack] strict"
#   String at the first line of our JavaScript file.
#
# The target layout moves the hidden "_GoBack" bookmark out of the
# "use strict" sentence and into a brand-new, otherwise-empty paragraph
# appended at the very end of the document (after the
# "String at the first line..." paragraph).

# Step 1: remove the "_GoBack" bookmark from its current location
# (it sits between "use " and "strict").  It is a hidden bookmark so it
# does not show up in the Bookmarks collection by index, but it can
# still be referenced directly by name.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Step 2: add a new, empty paragraph right after the
# "String at the first line of our JavaScript file." paragraph by
# appending a paragraph mark to its text via Find/Replace (this keeps
# the new paragraph free of any stray runs, matching a genuine
# paragraph-mark-only paragraph).
$d.Content.Find.Execute("JavaScript file.", $true, $false, $false, $false, $false, $true, 1, $false, "JavaScript file.^p", 2)

# Step 3: re-create the "_GoBack" bookmark inside that new, final
# (empty) paragraph.
$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastParagraph.Range)
